$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1, matching the existing header formatting
# (bold font, thin border, centered horizontally, top vertical alignment)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$headerRange = $ws.Range("F1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# New boolean "outlier" flag columns F:H for data rows 2-9
$outliers = @{
    2 = @($false, $false, $false)
    3 = @($true,  $false, $false)
    4 = @($false, $false, $false)
    5 = @($false, $false, $false)
    6 = @($false, $false, $false)
    7 = @($true,  $false, $false)
    8 = @($false, $false, $false)
    9 = @($false, $false, $false)
}

foreach ($row in $outliers.Keys) {
    $vals = $outliers[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}
